# Append: 2026-01-03 12:36 JST
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-03 12:36:02"

# The scraper re-stamps the fetch time on every existing row each run.
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(5, 1).Value = $newTimestamp

# Append the newly scraped case as row 6.
$row = 6

$ws.Cells.Item($row, 1).Value = $newTimestamp
$ws.Cells.Item($row, 2).Value = "ビジネスマッチングサイト構築・運用の依頼"
$ws.Cells.Item($row, 3).Value = "システム開発"
$ws.Cells.Item($row, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item($row, 5).Value = "期限情報なし"

$ws.Hyperlinks.Add($ws.Cells.Item($row, 6), "https://www.lancers.jp/work/detail/5464329")
$ws.Cells.Item($row, 6).Style = "Hyperlink"

$ws.Cells.Item($row, 7).Value = 45
$ws.Cells.Item($row, 8).Value = "◇サイト"
